$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new weekly record: shift existing rows 39-91 down to 40-92, ---
# --- and populate new row 39 with the latest week's data.                 ---

# Update cells in rows 39-91 so each row now holds the data previously held
# by the row above it (row 39 gets a new date; the rest cascade downward).

# Row 39
$ws.Range("D39").Value = 44902

# Row 40
$ws.Range("D40").Value = 44897
$ws.Range("H40").Value = 'Sin especificar'
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 1100
$ws.Range("L40").Value = 1100
$ws.Range("M40").Value = 1100
$ws.Range("P40").Value = 1100

# Row 41
$ws.Range("D41").Value = 44489
$ws.Range("H41").Value = 'Verde'
$ws.Range("J41").Value = 4000
$ws.Range("K41").Value = 900
$ws.Range("L41").Value = 900
$ws.Range("M41").Value = 900
$ws.Range("P41").Value = 900

# Row 42
$ws.Range("D42").Value = 44876

# Row 43
$ws.Range("D43").Value = 44879
$ws.Range("J43").Value = 3000
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 1000
$ws.Range("M43").Value = 1000
$ws.Range("P43").Value = 1000

# Row 44
$ws.Range("J44").Value = 2000

# Row 45
$ws.Range("D45").Value = 44895
$ws.Range("K45").Value = 1100
$ws.Range("L45").Value = 1100
$ws.Range("M45").Value = 1100
$ws.Range("P45").Value = 1100

# Row 46
$ws.Range("D46").Value = 44855
$ws.Range("H46").Value = 'Sin especificar'
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 1200
$ws.Range("M46").Value = 1200
$ws.Range("O46").Value = 'Provincia de Linares'
$ws.Range("P46").Value = 1200

# Row 47
$ws.Range("D47").Value = 44516
$ws.Range("J47").Value = 3000
$ws.Range("K47").Value = 1000
$ws.Range("L47").Value = 1000
$ws.Range("M47").Value = 1000
$ws.Range("O47").Value = 'Provincia de Limarí'
$ws.Range("P47").Value = 1000

# Row 48
$ws.Range("D48").Value = 44467
$ws.Range("H48").Value = 'Verde'
$ws.Range("J48").Value = 2000
$ws.Range("K48").Value = 1800
$ws.Range("L48").Value = 1800
$ws.Range("M48").Value = 1800
$ws.Range("P48").Value = 1800

# Row 49
$ws.Range("D49").Value = 44831

# Row 50
$ws.Range("D50").Value = 44827
$ws.Range("H50").Value = 'Sin especificar'
$ws.Range("J50").Value = 1000
$ws.Range("K50").Value = 2000
$ws.Range("L50").Value = 2000
$ws.Range("M50").Value = 2000
$ws.Range("P50").Value = 2000

# Row 51
$ws.Range("D51").Value = 44509
$ws.Range("J51").Value = 6000
$ws.Range("K51").Value = 800
$ws.Range("L51").Value = 800
$ws.Range("M51").Value = 800
$ws.Range("P51").Value = 800

# Row 52
$ws.Range("D52").Value = 44522
$ws.Range("J52").Value = 3000
$ws.Range("K52").Value = 1200
$ws.Range("L52").Value = 1200
$ws.Range("M52").Value = 1200
$ws.Range("P52").Value = 1200

# Row 53
$ws.Range("D53").Value = 44176
$ws.Range("J53").Value = 2000
$ws.Range("K53").Value = 900
$ws.Range("L53").Value = 900
$ws.Range("M53").Value = 900
$ws.Range("P53").Value = 900

# Row 54
$ws.Range("D54").Value = 44461
$ws.Range("H54").Value = 'Verde'
$ws.Range("J54").Value = 2500
$ws.Range("K54").Value = 2000
$ws.Range("L54").Value = 2000
$ws.Range("M54").Value = 2000
$ws.Range("P54").Value = 2000

# Row 55
$ws.Range("D55").Value = 44873
$ws.Range("H55").Value = 'Sin especificar'
$ws.Range("J55").Value = 3000
$ws.Range("K55").Value = 1100
$ws.Range("L55").Value = 1100
$ws.Range("M55").Value = 1100
$ws.Range("P55").Value = 1100

# Row 56
$ws.Range("D56").Value = 44491
$ws.Range("K56").Value = 850
$ws.Range("L56").Value = 850
$ws.Range("M56").Value = 850
$ws.Range("P56").Value = 850

# Row 57
$ws.Range("D57").Value = 44476
$ws.Range("H57").Value = 'Verde'
$ws.Range("J57").Value = 5000
$ws.Range("K57").Value = 1000
$ws.Range("L57").Value = 1100
$ws.Range("M57").Value = 1040
$ws.Range("P57").Value = 1040

# Row 58
$ws.Range("D58").Value = 44837
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 1600
$ws.Range("L58").Value = 1600
$ws.Range("M58").Value = 1600
$ws.Range("O58").Value = 'Provincia de Linares'
$ws.Range("P58").Value = 1600

# Row 59
$ws.Range("D59").Value = 44860
$ws.Range("J59").Value = 4000
$ws.Range("K59").Value = 1000
$ws.Range("L59").Value = 1000
$ws.Range("M59").Value = 1000
$ws.Range("O59").Value = 'Región del Maule'
$ws.Range("P59").Value = 1000

# Row 60
$ws.Range("D60").Value = 44832
$ws.Range("H60").Value = 'Sin especificar'
$ws.Range("J60").Value = 1500
$ws.Range("K60").Value = 2000
$ws.Range("L60").Value = 2000
$ws.Range("M60").Value = 2000
$ws.Range("P60").Value = 2000

# Row 61
$ws.Range("D61").Value = 44477
$ws.Range("H61").Value = 'Verde'
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 1000
$ws.Range("L61").Value = 1000
$ws.Range("M61").Value = 1000
$ws.Range("P61").Value = 1000

# Row 62
$ws.Range("D62").Value = 44894
$ws.Range("K62").Value = 1100
$ws.Range("L62").Value = 1100
$ws.Range("M62").Value = 1100
$ws.Range("N62").Value = '$/kilo'
$ws.Range("P62").Value = 1100

# Row 63
$ws.Range("D63").Value = 44883
$ws.Range("H63").Value = 'Sin especificar'
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 1000
$ws.Range("L63").Value = 1000
$ws.Range("M63").Value = 1000
$ws.Range("N63").Value = '$/atado'
$ws.Range("P63").Value = 1000

# Row 64
$ws.Range("D64").Value = 44497
$ws.Range("H64").Value = 'Verde'
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 900
$ws.Range("L64").Value = 900
$ws.Range("M64").Value = 900
$ws.Range("O64").Value = 'Provincia de Linares'
$ws.Range("P64").Value = 900

# Row 65
$ws.Range("D65").Value = 44882
$ws.Range("H65").Value = 'Sin especificar'
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 1000
$ws.Range("L65").Value = 1000
$ws.Range("M65").Value = 1000
$ws.Range("N65").Value = '$/kilo'
$ws.Range("O65").Value = 'Región del Maule'
$ws.Range("P65").Value = 1000
$ws.Range("Q65").Value = 1

# Row 66
$ws.Range("D66").Value = 44481
$ws.Range("N66").Value = '$/caja 10 kilos'
$ws.Range("P66").Value = 90
$ws.Range("Q66").Value = 10

# Row 67
$ws.Range("D67").Value = 44496
$ws.Range("J67").Value = 4000

# Row 68
$ws.Range("D68").Value = 44498
$ws.Range("H68").Value = 'Verde'
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 900
$ws.Range("L68").Value = 900
$ws.Range("M68").Value = 900
$ws.Range("P68").Value = 900

# Row 69
$ws.Range("D69").Value = 44875
$ws.Range("J69").Value = 3000
$ws.Range("K69").Value = 1000
$ws.Range("L69").Value = 1000
$ws.Range("M69").Value = 1000
$ws.Range("P69").Value = 1000

# Row 70
$ws.Range("D70").Value = 44839
$ws.Range("J70").Value = 2500
$ws.Range("K70").Value = 1300
$ws.Range("L70").Value = 1300
$ws.Range("M70").Value = 1300
$ws.Range("P70").Value = 1300

# Row 71
$ws.Range("D71").Value = 44847
$ws.Range("H71").Value = 'Sin especificar'
$ws.Range("J71").Value = 3000

# Row 72
$ws.Range("D72").Value = 44474
$ws.Range("J72").Value = 5000

# Row 73
$ws.Range("D73").Value = 44515
$ws.Range("K73").Value = 1200
$ws.Range("L73").Value = 1200
$ws.Range("M73").Value = 1200
$ws.Range("O73").Value = 'Provincia de Linares'
$ws.Range("P73").Value = 1200

# Row 74
$ws.Range("D74").Value = 44523
$ws.Range("K74").Value = 1100
$ws.Range("L74").Value = 1100
$ws.Range("M74").Value = 1100
$ws.Range("O74").Value = 'Región del Maule'
$ws.Range("P74").Value = 1100

# Row 75
$ws.Range("D75").Value = 44168
$ws.Range("J75").Value = 3000
$ws.Range("L75").Value = 1000
$ws.Range("M75").Value = 1000
$ws.Range("P75").Value = 1000

# Row 76
$ws.Range("D76").Value = 44475
$ws.Range("J76").Value = 5000
$ws.Range("K76").Value = 1000
$ws.Range("L76").Value = 1100
$ws.Range("M76").Value = 1040
$ws.Range("O76").Value = 'Provincia de Linares'
$ws.Range("P76").Value = 1040

# Row 77
$ws.Range("D77").Value = 44490
$ws.Range("H77").Value = 'Verde'
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 850
$ws.Range("L77").Value = 900
$ws.Range("M77").Value = 875
$ws.Range("O77").Value = 'Región del Maule'
$ws.Range("P77").Value = 875

# Row 78
$ws.Range("I78").Value = 'Primera'
$ws.Range("J78").Value = 2500
$ws.Range("K78").Value = 1100
$ws.Range("L78").Value = 1100
$ws.Range("M78").Value = 1100
$ws.Range("P78").Value = 1100

# Row 79
$ws.Range("D79").Value = 44868
$ws.Range("I79").Value = 'Segunda'
$ws.Range("J79").Value = 2000
$ws.Range("K79").Value = 1000
$ws.Range("L79").Value = 1000
$ws.Range("M79").Value = 1000
$ws.Range("P79").Value = 1000

# Row 80
$ws.Range("D80").Value = 44830
$ws.Range("H80").Value = 'Sin especificar'
$ws.Range("J80").Value = 500
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 2000
$ws.Range("M80").Value = 2000
$ws.Range("O80").Value = 'Provincia de Linares'
$ws.Range("P80").Value = 2000

# Row 81
$ws.Range("D81").Value = 44512
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 800
$ws.Range("L81").Value = 800
$ws.Range("M81").Value = 800
$ws.Range("O81").Value = 'Región del Maule'
$ws.Range("P81").Value = 800

# Row 82
$ws.Range("D82").Value = 44517
$ws.Range("J82").Value = 4000
$ws.Range("K82").Value = 1100
$ws.Range("L82").Value = 1100
$ws.Range("M82").Value = 1100
$ws.Range("O82").Value = 'Provincia de Linares'
$ws.Range("P82").Value = 1100

# Row 83
$ws.Range("D83").Value = 44172
$ws.Range("J83").Value = 2000
$ws.Range("K83").Value = 1000
$ws.Range("M83").Value = 1000
$ws.Range("O83").Value = 'Región del Maule'
$ws.Range("P83").Value = 1000

# Row 84
$ws.Range("D84").Value = 44482
$ws.Range("J84").Value = 4000
$ws.Range("K84").Value = 900
$ws.Range("M84").Value = 950
$ws.Range("O84").Value = 'Provincia de Linares'
$ws.Range("P84").Value = 950

# Row 85
$ws.Range("D85").Value = 44524
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 1000
$ws.Range("L85").Value = 1000
$ws.Range("M85").Value = 1000
$ws.Range("O85").Value = 'Región del Maule'
$ws.Range("P85").Value = 1000

# Row 86
$ws.Range("D86").Value = 44505
$ws.Range("H86").Value = 'Verde'
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 800
$ws.Range("L86").Value = 800
$ws.Range("M86").Value = 800
$ws.Range("P86").Value = 800

# Row 87
$ws.Range("D87").Value = 44881
$ws.Range("H87").Value = 'Sin especificar'
$ws.Range("K87").Value = 1000
$ws.Range("L87").Value = 1000
$ws.Range("M87").Value = 1000
$ws.Range("P87").Value = 1000

# Row 88
$ws.Range("D88").Value = 44468
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 1500
$ws.Range("M88").Value = 1500
$ws.Range("P88").Value = 1500

# Row 89
$ws.Range("D89").Value = 44161
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 1000
$ws.Range("L89").Value = 1000
$ws.Range("M89").Value = 1000
$ws.Range("P89").Value = 1000

# Row 90
$ws.Range("D90").Value = 44460
$ws.Range("H90").Value = 'Verde'
$ws.Range("J90").Value = 2000
$ws.Range("K90").Value = 2000
$ws.Range("L90").Value = 2000
$ws.Range("M90").Value = 2000
$ws.Range("P90").Value = 2000

# Row 91
$ws.Range("D91").Value = 44848
$ws.Range("H91").Value = 'Sin especificar'
$ws.Range("J91").Value = 3000
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 1000
$ws.Range("M91").Value = 1000
$ws.Range("O91").Value = 'Provincia de Linares'
$ws.Range("P91").Value = 1000

# Re-apply the date number format on the shifted date cells so the new row
# retains the same display format as the rest of column D.
$ws.Range("D39:D91").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 92: the record that used to live in row 91 before the insert.
$ws.Range("A92").Value = 5
$ws.Range("B92").Value = 'Macroferia Regional de Talca'
$ws.Range("C92").Value = 'Maule'
$ws.Range("D92").Value = 44519
$ws.Range("E92").Value = 7
$ws.Range("F92").Value = 300000000
$ws.Range("G92").Value = 'Espárragos'
$ws.Range("H92").Value = 'Verde'
$ws.Range("I92").Value = 'Primera'
$ws.Range("J92").Value = 4000
$ws.Range("K92").Value = 1100
$ws.Range("L92").Value = 1100
$ws.Range("M92").Value = 1100
$ws.Range("N92").Value = '$/kilo'
$ws.Range("O92").Value = 'Región del Maule'
$ws.Range("P92").Value = 1100
$ws.Range("Q92").Value = 1
$ws.Range("R92").Value = 'Hortaliza'

# Row 92's date cell needs the same date format as the rest of column D.
$ws.Range("D92").NumberFormat = "YYYY-MM-DD HH:MM:SS"
